$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = 78

$ws.Cells.Item($r, 1).Value = 2
$ws.Cells.Item($r, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item($r, 3).Value = "Coquimbo"
$ws.Cells.Item($r, 4).Value = "2022-09-28"
$ws.Cells.Item($r, 4).NumberFormat = $ws.Cells.Item(77, 4).NumberFormat
$ws.Cells.Item($r, 5).Value = 4
$ws.Cells.Item($r, 6).Value = 100112022
$ws.Cells.Item($r, 7).Value = "Arveja Verde"
$ws.Cells.Item($r, 8).Value = "Perfection"
$ws.Cells.Item($r, 9).Value = "Primera"
$ws.Cells.Item($r, 10).Value = 600
$ws.Cells.Item($r, 11).Value = 23000
$ws.Cells.Item($r, 12).Value = 25000
$ws.Cells.Item($r, 13).Value = 24000
$ws.Cells.Item($r, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item($r, 15).Value = "Provincia de Limarí"
$ws.Cells.Item($r, 16).Value = 960
$ws.Cells.Item($r, 17).Value = 25
$ws.Cells.Item($r, 18).Value = "Hortaliza"
